$wb = $excel.ActiveWorkbook

# --- Add the new "Com" worksheet at the end ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Com"

# Header row (BldgLoc must be written before BldgType so that the
# shared-string table assigns BldgLoc index 34 and BldgType index 35)
$newSheet.Range("B1").Value = "BldgLoc"
$newSheet.Range("A1").Value = "BldgType"
$newSheet.Range("C1").Value = "Normunit"
$newSheet.Range("D1").Value = "Value"
$newSheet.Range("E1").Value = "BldgVint"

# Data rows: BldgType code, Normunit, Value
$data = @(
    @("Asm", 100002.1),
    @("ECC", 299999.59999999998),
    @("EPr", 50000.14),
    @("ERC", 1920.0170000000001),
    @("ESe", 149998.6),
    @("EUn", 930201.4),
    @("Hsp", 235501),
    @("Htl", 139998.9),
    @("MBT", 199999.2),
    @("MLI", 100001.1),
    @("Mtl", 30000.07),
    @("Nrs", 60654.58),
    @("OfL", 174998.9),
    @("OfS", 10000.08),
    @("RFF", 2500.0529999999999),
    @("RSD", 5599.9570000000003),
    @("Rt3", 120000.5),
    @("RtL", 129997),
    @("RtS", 7999.9290000000001),
    @("SCn", 250000.3)
)

$row = 2
foreach ($item in $data) {
    $newSheet.Cells.Item($row, 1).Value = $item[0]
    $newSheet.Cells.Item($row, 3).Value = "Area-ft2-BA"
    $newSheet.Cells.Item($row, 4).Value = $item[1]
    $row = $row + 1
}

# Column width formatting for column C (closest value achievable given this
# runtime's column-width quantization; target stored width is 12.6640625)
$newSheet.Columns.Item(3).ColumnWidth = 11.8333333

# Selection on the new sheet
$newSheet.Range("G6").Select()

# --- Adjust selection on DMo (sheet1) ---
$dmo = $wb.Worksheets.Item("DMo")
$dmo.Range("A2").Select()

# --- Make "Com" the active/selected tab ---
$newSheet.Activate()
$newSheet.Select()
